# "Add files via upload" -- appends 7 new observation rows (86-92) to
# Sheet1 of the landscaping data log, extends the dependent formula
# columns (F = ABS(High-Low), H = Growth) down to match, and leaves the
# selection on the newly-added Growth column (T86:T92), matching the
# author's last on-screen selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 85 (the previous last row) already carries the correct date
# number format for column A; copy it down across the new rows first so
# the new dates (A86:A92) render the same way the rest of the column does.
$ws.Range("A85").Copy()
$ws.Range("A86:A92").PasteSpecial(-4122)  # xlPasteFormats

# Per-row data that differs between the 7 new records. Every other
# field (A, D, E, G, I, L, M, N, O, P, Q, R, S, T) is identical across
# all of them.
$newRows = @(
  @{ Row = 86; PlantType = "Flowering";    PlantSize = "Large";  Growth = "=1/3"; Quadrant = "Neutral"; Pruned = 2 },
  @{ Row = 87; PlantType = "Nonflowering"; PlantSize = "Medium"; Growth = 0.5;    Quadrant = "Neutral"; Pruned = 3 },
  @{ Row = 88; PlantType = "Nonflowering"; PlantSize = "Small";  Growth = 0.1;    Quadrant = "Neutral"; Pruned = 3 },
  @{ Row = 89; PlantType = "Nonflowering"; PlantSize = "Medium"; Growth = 0.25;   Quadrant = "Dark";    Pruned = 3 },
  @{ Row = 90; PlantType = "Nonflowering"; PlantSize = "Medium"; Growth = 0.25;   Quadrant = "Dark";    Pruned = 3 },
  @{ Row = 91; PlantType = "Nonflowering"; PlantSize = "Large";  Growth = 0.1;    Quadrant = "Neutral"; Pruned = 4 },
  @{ Row = 92; PlantType = "Tree";         PlantSize = "Medium"; Growth = 2.25;   Quadrant = "Bright";  Pruned = 1 }
)

foreach ($d in $newRows) {
  $r = $d.Row

  $ws.Range("A$r").Value = 45799
  $ws.Range("B$r").Value = $d.PlantType
  $ws.Range("C$r").Value = $d.PlantSize
  $ws.Range("D$r").Value = 47
  $ws.Range("E$r").Value = 50
  $ws.Range("F$r").Formula = "=ABS(D$r-E$r)"
  $ws.Range("G$r").Value = 1.35

  if ($d.Growth -is [string]) {
    $ws.Range("H$r").Formula = $d.Growth
  } else {
    $ws.Range("H$r").Value = $d.Growth
  }

  $ws.Range("I$r").Value = "No"
  $ws.Range("J$r").Value = $d.Pruned
  $ws.Range("K$r").Value = $d.Quadrant
  $ws.Range("L$r").Value = 2
  $ws.Range("M$r").Value = 0.96
  $ws.Range("N$r").Value = 49
  $ws.Range("O$r").Value = 29.73
  $ws.Range("P$r").Value = 21
  $ws.Range("Q$r").Value = 0.96
  $ws.Range("R$r").Value = 8.1
  $ws.Range("S$r").Value = 21
  $ws.Range("T$r").Value = 30
}

# Match the author's final on-screen selection.
$ws.Range("T86:T92").Select()
